# Applies the "Automatic update of files." change:
#  - Column C (Förändrad) on rows 2-12 changes from 46070 to 46072
#  - Rows 5-9 (Beteckning/Datum/Area) get re-ordered to a new arrangement,
#    while keeping each record's own Beteckning (A), Datum (B) and Area (G)
#    values bundled together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Förändrad (C) value for every data row
$nyttDatum = 46072

# Update column C for all data rows (2 through 12)
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 3).Value2 = $nyttDatum
}

# New arrangement for rows 5-9 : Beteckning (A), Datum (B), Area (G)
# Row 5 <- old row 6 (A 50762-2025)
$ws.Cells.Item(5, 1).Value2 = "A 50762-2025"
$ws.Cells.Item(5, 2).Value2 = 45946
$ws.Cells.Item(5, 7).Value2 = 2.7

# Row 6 <- old row 5 (A 14517-2023)
$ws.Cells.Item(6, 1).Value2 = "A 14517-2023"
$ws.Cells.Item(6, 2).Value2 = 45012
$ws.Cells.Item(6, 7).Value2 = 0.6

# Row 7 <- old row 9 (A 14516-2023)
$ws.Cells.Item(7, 1).Value2 = "A 14516-2023"
$ws.Cells.Item(7, 2).Value2 = 45012.86600694444
$ws.Cells.Item(7, 7).Value2 = 0.4

# Row 8 <- old row 7 (A 8679-2026)
$ws.Cells.Item(8, 1).Value2 = "A 8679-2026"
$ws.Cells.Item(8, 2).Value2 = 46066
$ws.Cells.Item(8, 7).Value2 = 2.1

# Row 9 <- old row 8 (A 8929-2026)
$ws.Cells.Item(9, 1).Value2 = "A 8929-2026"
$ws.Cells.Item(9, 2).Value2 = 46069.34543981482
$ws.Cells.Item(9, 7).Value2 = 1.2
